$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'51.856.53"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.21%  "
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').Value = "'2.781.90"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -1.93%  "
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('E4').Value = "'  -0.05%  "
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').Value = "'358.29"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +1.34%  "
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').Value = "'109.20"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -3.82%  "
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('D7').Value = "'0.564"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.16%  "
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('E8').Value = "'  -0.04%  "
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('D9').Value = "'0.596"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.72%  "
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('D10').Value = "'40.05"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -3.86%  "
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').Value = "'0.0852"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.00%  "
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('E12').Value = "'  +0.68%  "
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('D13').Value = "'19.50"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -2.78%  "
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('D14').Value = "'7.59"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -1.70%  "
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').Value = "'3.215.88"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -2.10%  "
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').Value = "'2.793.70"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -1.31%  "
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('E17').Value = "'  +4.02%  "
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').Value = "'51.777.15"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.31%  "
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('D19').Value = "'7.43"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.69%  "
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').Value = "'3.12"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -1.48%  "
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('D21').Value = "'13.02"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -3.51%  "
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('D22').Value = "'0.0₃0978"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -1.70%  "
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('D23').Value = "'274.41"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.31%  "
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('D24').Value = "'70.23"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.68%  "
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('E25').Value = "'  -1.62%  "
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('D26').Value = "'26.69"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.35%  "
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('E27').Value = "'  -0.14%  "
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('D28').Value = "'10.18"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -1.29%  "
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'2.23"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -1.10%  "
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = "'0.144"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +3.65%  "
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('D31').Value = "'0.0465"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +4.63%  "
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('D32').Value = "'51.40"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.52%  "
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('D33').Value = "'33.98"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.47%  "
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('D34').Value = "'5.69"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -2.36%  "
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('D35').Value = "'0.0844"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +2.23%  "
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('D36').Value = "'5.24"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +7.18%  "
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('E37').Value = "'  +0.03%  "
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('E38').Value = "'  +0.98%  "
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('D39').Value = "'18.10"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.04%  "
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('E40').Value = "'  -3.86%  "
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('D41').Value = "'2.56"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.54%  "
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('E42').Value = "'  -1.15%  "
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('E43').Value = "'  -2.50%  "
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('D44').Value = "'121.73"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -3.12%  "
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('D45').Value = "'21.93"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -7.09%  "
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('D46').Value = "'2.069.57"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.73%  "
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('D47').Value = "'3.24"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -2.55%  "
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('D48').Value = "'2.18"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -5.75%  "
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('D49').Value = "'5.69"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.07%  "
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('E50').Value = "'  -0.73%  "
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('D51').Value = "'8.94"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.07%  "
$ws.Range('E51').Style = 'Normal'
